$d = $word.ActiveDocument

# Replace-SubtitleText finds every occurrence of $old in the document
# and sets it to $new. Assigning directly to the found Range's .Text
# (instead of passing $new as Find.Execute's Replace argument) avoids
# Word's smart-quote AutoFormat mangling straight apostrophes.
function Replace-SubtitleText($old, $new, $wholeWord) {
    $keepGoing = $true
    while ($keepGoing) {
        $rng = $d.Content
        $rng.Find.Execute($old, $true, $false, $false, $wholeWord, $false, $true, 1, $false) | Out-Null
        if ($rng.Find.Found) {
            $rng.Text = $new
        } else {
            $keepGoing = $false
        }
    }
}

Replace-SubtitleText "Tatizo la mchwa - manukuu:" "The ants problem - subtitles:" $false
Replace-SubtitleText "Mazungumzo huanza kwa sekunde 40 kwa hivyo niliongeza sekunde 27 kwa nyakati kama zilivyokuwa - John Argentino" "The dialogue starts at 40 seconds in so I added 27 seconds to the times as they were - John Argentino" $false
Replace-SubtitleText "[Muziki]" "[Music]" $false
Replace-SubtitleText "sawa kwa hivyo mafumbo nitaenda" "okay so the puzzles I'm going to" $false
Replace-SubtitleText "changamoto uliyonayo ni mbili za msingi" "challenge you with are two basic" $false
Replace-SubtitleText "matoleo ya fumbo ngumu zaidi" "versions of a more complicated puzzle" $false
Replace-SubtitleText "inayojulikana kama fumbo la mchwa, ambalo mimi ni" "known as the ants puzzle, which I'm" $false
Replace-SubtitleText "pengine kwenda kujadili katika tofauti" "probably going to discuss in a different" $false
Replace-SubtitleText "video. Ngoja nimalizie kuandika" "video. Let me just finish writing down" $false
Replace-SubtitleText "kichwa na, vizuri, naweza hata kuchora a" "the title and, well, I can even draw a" $false
Replace-SubtitleText "mchwa mdogo hapa. sawa, tupate" "little ant right here. okay, let's get" $false
Replace-SubtitleText "imeanza! Kama nilivyosema nitajadili" "started! As I said I'm going to discuss" $false
Replace-SubtitleText "mafumbo mawili katika fumbo la kwanza hapo" "two puzzles in the first puzzle there" $false
Replace-SubtitleText "ni mchwa wawili kwenye kinyesi cha juu sana: aina" "are two ants on a very high stool: a sort" $false
Replace-SubtitleText "ya Mlima, gorofa juu na mbili" "of Mountain, flat on the top with two" $false
Replace-SubtitleText "miamba mikali kwa pande zote mbili. Gorofa" "steep cliffs to both the sides. The flat" $false
Replace-SubtitleText "kilele ni mita moja upana wa mchwa wawili hoja" "peak is one meter wide the two ants move" $false
Replace-SubtitleText "kwa kasi, tuiite V, ambayo ni" "with a velocity, let's call it V, which is" $false
Replace-SubtitleText "sawa kwa wote wawili na hiyo ni" "the same for both of them and that is" $false
Replace-SubtitleText "sawa na sentimita moja kwa sekunde. Wewe" "equal to one centimeter per second. You" $false
Replace-SubtitleText "inaweza kuamua mwelekeo kuelekea kila mmoja" "can decide the direction towards each" $false
Replace-SubtitleText "mchwa husogea ikiwa ni kulia au kushoto na" "ant moves if it is right or left and" $false
Replace-SubtitleText "wapi hasa kuweka mchwa wawili kwenye" "where exactly to place the two ants on the" $false
Replace-SubtitleText "juu ya mlima. Kusudi lako ni" "top of the mountain. Your purpose is to" $false
Replace-SubtitleText "fanya wakati mchwa wa mwisho huchukua hapo awali" "make the time the last ant takes before" $false
Replace-SubtitleText "kuanguka kwa muda mrefu iwezekanavyo. Mchwa hawawezi" "falling the longest possible. Ants cannot" $false
Replace-SubtitleText "tulia: lazima wahamie kulia au" "be still: they must move to the right or" $false
Replace-SubtitleText "upande wa kushoto lakini lazima wasogee na baada" "to the left but they must move and after" $false
Replace-SubtitleText "wakikutana wanageuka na" "meeting each other they turn around and" $false
Replace-SubtitleText "endelea kusonga na sawa lakini kinyume" "keep moving with the same but opposite" $false
Replace-SubtitleText "kwa hivyo tena ni nafasi gani sahihi" "so again what are the precise positions" $false
Replace-SubtitleText "ambapo ninapaswa kuwaweka mchwa wawili ndani" "where I should place the two ants in" $false
Replace-SubtitleText "ili kupata muda mrefu zaidi kabla ya" "order to get the longest time before the" $false
Replace-SubtitleText "chungu mwisho huanguka? Fumbo la pili ni" "last ant falls? The second puzzle is" $false
Replace-SubtitleText "kimsingi ni sawa lakini sasa tuna tatu" "basically the same but now we have three" $false
Replace-SubtitleText "mchwa badala ya wawili." "ants instead of two." $false
Replace-SubtitleText "Kama kabla ya mchwa kasi ni moja" "As before the ants velocity is one" $false
Replace-SubtitleText "sentimita kwa sekunde, kila mchwa hugeuka" "centimeter per second, every ant turns" $false
Replace-SubtitleText "karibu baada ya kukutana na mchwa mwingine na" "around after meeting another ant and" $false
Replace-SubtitleText "kilele kina upana wa mita moja. Hivyo, ni nini" "the peak is one meter wide. So, what are" $false
Replace-SubtitleText "sasa nafasi sahihi" "now the precise positions" $false
Replace-SubtitleText "Ninapaswa kuweka mchwa watatu kwa mpangilio" "I should place the three ants in order" $false
Replace-SubtitleText "kupata muda mrefu zaidi kabla ya mwisho" "to get the longest time before the last" $false
Replace-SubtitleText "chungu huanguka chini? Natumaini ulifurahia hili" "ant falls down? I hope you enjoyed this" $false
Replace-SubtitleText "video fanya bora na bahati nzuri" "video do your best and good luck" $false
Replace-SubtitleText "kasi" "velocity" $true
